$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New diary entries (rows 63-66) appended to the "Summer 2024 Working Hours" log.
# Column A = date, Column B = day-of-week letter, Column C = hours, Column E = notes.

$ws.Range("A63").Value = 45504
$ws.Range("A63").NumberFormat = "d-mmm"
$ws.Range("B63").Value = "W"
$ws.Range("B63").NumberFormat = "d-mmm"
$ws.Range("C63").Value = 4
$ws.Range("E63").Value = "Building datasets for COSC 221, DATA 301 & 311. Used STAT 230 too. Built grid search to find patterns between best hyperparameters. Started code"
$ws.Range("E63").WrapText = $true
$ws.Rows.Item(63).RowHeight = 42.75

$ws.Range("A64").Value = 45505
$ws.Range("A64").NumberFormat = "d-mmm"
$ws.Range("B64").Value = "T"
$ws.Range("B64").NumberFormat = "d-mmm"
$ws.Range("C64").Value = 4
$ws.Range("E64").Value = "A lower shrinkage rate and slightly higher interacton depth performed better,  Minimum number of observations in terminal nodes (n.minobsinnode) seemed to move around with no pattern. Repeated for STAT 303 and MATH 221"
$ws.Range("E64").WrapText = $true
$ws.Rows.Item(64).RowHeight = 71.25

$ws.Range("A65").Value = 45507
$ws.Range("A65").NumberFormat = "d-mmm"
$ws.Range("B65").Value = "S"
$ws.Range("B65").NumberFormat = "d-mmm"
$ws.Range("C65").Value = 4
$ws.Range("E65").Value = "MATH 221 didn't seem to perform well depite having the most students having taken the course. Doing separate testing on on MATH 221 on the main function to see if the setup isn't working."
$ws.Range("E65").WrapText = $true
$ws.Rows.Item(65).RowHeight = 57

$ws.Range("A66").Value = 45509
$ws.Range("A66").NumberFormat = "d-mmm"
$ws.Range("B66").Value = "M"
$ws.Range("B66").NumberFormat = "d-mmm"
$ws.Range("C66").Value = 4
$ws.Range("E66").Value = "Re running  grid search on all datasets and re evaluating RMSEs. The RMSEs dropped 0.5-1 for each of the courses and similar observations with hyperparameters were found."
$ws.Range("E66").WrapText = $true
$ws.Rows.Item(66).RowHeight = 57

# Update the view so the newly added rows are visible, matching the author's
# last scroll position / selection when they saved the file.
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Range("E65").Select()
